$wb = $excel.ActiveWorkbook

# --- "Inactive" sheet: add a new todo item (Id 30) at the top of the list ---
$ws = $wb.Worksheets.Item("Inactive")

# Duplicate the existing row 2 (same Status/Category/Created/Done values &
# formatting as the rest of the inactive items) and insert it above, shifting
# the old rows 2-3 down to 3-4.
$ws.Range("A2:F2").Copy()
$ws.Range("A2:F2").Insert()
$ws.Range("A2:F2").ClearFormats()

# Fill in the new item's own Id and Title; Status/Category/Created/Done stay
# the same as every other inactive (done) task.
$ws.Cells.Item(2, 1).Value = 30
$ws.Cells.Item(2, 2).Value = "resize image to fit window as it resizes"

# --- "Config" sheet: bump Max Id from 29 to 30 to reflect the new task ---
$cfg = $wb.Worksheets.Item("Config")
$cfg.Cells.Item(2, 6).Value = 30
